$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $val) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $val
    $cellRange.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "68.439.27"
$ws.Range("E2").Value = "  -1.89%  "

Set-TextValue $ws.Range("D3") "2.450.99"
$ws.Range("E3").Value = "  -2.25%  "

$ws.Range("E4").Value = "  +0.00%  "

Set-TextValue $ws.Range("D5") "564.05"
$ws.Range("E5").Value = "  -2.09%  "

Set-TextValue $ws.Range("D6") "163.97"
$ws.Range("E6").Value = "  -2.05%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  -1.37%  "

$ws.Range("E9").Value = "  -6.41%  "

$ws.Range("E10").Value = "  -2.04%  "

Set-TextValue $ws.Range("D11") "0.341"
$ws.Range("E11").Value = "  -4.44%  "

$ws.Range("E12").Value = "  -2.69%  "

Set-TextValue $ws.Range("D13") "2.902.92"
$ws.Range("E13").Value = "  -2.17%  "

Set-TextValue $ws.Range("D14") "68.404.19"
$ws.Range("E14").Value = "  -1.65%  "

$ws.Range("E15").Value = "  -4.20%  "

$ws.Range("E16").Value = "  -5.17%  "

Set-TextValue $ws.Range("D17") "2.458.52"
$ws.Range("E17").Value = "  -1.25%  "

$ws.Range("E18").Value = "  -2.23%  "

Set-TextValue $ws.Range("D19") "344.58"
$ws.Range("E19").Value = "  -1.42%  "

Set-TextValue $ws.Range("D20") "7.15"
$ws.Range("E20").Value = "  -4.60%  "

$ws.Range("E21").Value = "  -2.13%  "

$ws.Range("E22").Value = "  -3.10%  "

Set-TextValue $ws.Range("D23") "1.00"
$ws.Range("E23").Value = "  +0.00%  "

Set-TextValue $ws.Range("D24") "68.10"
$ws.Range("E24").Value = "  -3.44%  "

Set-TextValue $ws.Range("D25") "3.75"
$ws.Range("E25").Value = "  -5.20%  "

Set-TextValue $ws.Range("D26") "2.580.14"
$ws.Range("E26").Value = "  +1.13%  "

$ws.Range("E27").Value = "  +3.56%  "

Set-TextValue $ws.Range("D28") "8.26"
$ws.Range("E28").Value = "  -6.23%  "

$ws.Range("E29").Value = "  -5.73%  "

Set-TextValue $ws.Range("D30") "7.31"
$ws.Range("E30").Value = "  -6.86%  "

Set-TextValue $ws.Range("D31") "438.14"
$ws.Range("E31").Value = "  -4.89%  "

$ws.Range("E32").Value = "  -3.23%  "

Set-TextValue $ws.Range("D33") "0.999"
$ws.Range("E33").Value = "  -0.03%  "

Set-TextValue $ws.Range("D34") "1.69"
$ws.Range("E34").Value = "  -2.79%  "

Set-TextValue $ws.Range("D35") "3.05"
$ws.Range("E35").Value = "  +104.33%  "

Set-TextValue $ws.Range("D36") "156.60"

$ws.Range("E39").Value = "  -5.73%  "

$ws.Range("E40").Value = "  -3.21%  "

$ws.Range("E41").Value = "  -3.78%  "

$ws.Range("E42").Value = "  -4.34%  "

$ws.Range("E43").Value = "  -4.31%  "

$ws.Range("E44").Value = "  +5.36%  "

$ws.Range("E45").Value = "  -5.20%  "

Set-TextValue $ws.Range("D46") "135.30"
$ws.Range("E46").Value = "  -4.40%  "

$ws.Range("E47").Value = "  -3.18%  "

$ws.Range("E48").Value = "  -2.48%  "

$ws.Range("E49").Value = "  -6.34%  "

$ws.Range("E50").Value = "  -2.54%  "

$ws.Range("E51").Value = "  -1.51%  "
